$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 2213.8
$ws.Range("I135").Value = 2213.8
$ws.Range("K135").Value = 19924.2
$ws.Range("M135").Value = -17389.2
$ws.Range("H138").Value = 1229871.2
$ws.Range("I138").Value = 1647.5555
$ws.Range("J138").Value = 1672031.9
$ws.Range("K138").Value = 4942.666499999999
$ws.Range("L138").Value = 5016095.699999999
$ws.Range("M138").Value = 197.3335000000006
$ws.Range("N138").Value = -5026375.699999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3962.4614
$ws.Range("J63").Value = 4901
$ws.Range("L63").Value = 4901
$ws.Range("N63").Value = -6273
$ws.Range("H66").Value = 3962.4614
$ws.Range("J66").Value = 4901
$ws.Range("L66").Value = 24505
$ws.Range("N66").Value = -31369
$ws.Range("H74").Value = 885.75
$ws.Range("I74").Value = 457.5
$ws.Range("K74").Value = 457.5
$ws.Range("M74").Value = 416.5
$ws.Range("H77").Value = 885.75
$ws.Range("I77").Value = 457.5
$ws.Range("K77").Value = 2287.5
$ws.Range("M77").Value = 2080.5
$ws.Range("H122").Value = 914
$ws.Range("I122").Value = 914
$ws.Range("K122").Value = 2742
$ws.Range("M122").Value = -292
$ws.Range("H124").Value = 61960
$ws.Range("J124").Value = 61960
$ws.Range("L124").Value = 61960
$ws.Range("N124").Value = -71780

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H44").Value = 19999.5
$ws.Range("J44").Value = 19999.5
$ws.Range("L44").Value = 19999.5
$ws.Range("N44").Value = -20993.5
$ws.Range("H94").Value = 3600.0833
$ws.Range("I94").Value = 1020.2
$ws.Range("J94").Value = 16499.5
$ws.Range("K94").Value = 1020.2
$ws.Range("L94").Value = 16499.5
$ws.Range("M94").Value = -569.2
$ws.Range("N94").Value = -17401.5
$ws.Range("H107").Value = 3337.75
$ws.Range("I107").Value = 2896.5
$ws.Range("K107").Value = 2896.5
$ws.Range("M107").Value = -976.5
$ws.Range("H122").Value = 69140
$ws.Range("J122").Value = 69140
$ws.Range("L122").Value = 69140
$ws.Range("N122").Value = -78940
$ws.Range("H134").Value = 3781.182
$ws.Range("I134").Value = 2932.6667
$ws.Range("J134").Value = 4799.4
$ws.Range("K134").Value = 8798.000100000001
$ws.Range("L134").Value = 14398.2
$ws.Range("M134").Value = -6263.000100000001
$ws.Range("N134").Value = -19468.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2374.625
$ws.Range("I22").Value = 2149.75
$ws.Range("J22").Value = 2599.5
$ws.Range("K22").Value = 2149.75
$ws.Range("L22").Value = 2599.5
$ws.Range("M22").Value = -1799.75
$ws.Range("N22").Value = -3299.5
$ws.Range("H94").Value = 2090.9092
$ws.Range("I94").Value = 2532.7144
$ws.Range("K94").Value = 2532.7144
$ws.Range("M94").Value = -2081.7144
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H121").Value = 21249.5
$ws.Range("J121").Value = 21249.5
$ws.Range("L121").Value = 21249.5
$ws.Range("N121").Value = -23869.5
$ws.Range("H122").Value = 3241.5
$ws.Range("I122").Value = 2491.8572
$ws.Range("J122").Value = 3991.1428
$ws.Range("K122").Value = 7475.571599999999
$ws.Range("L122").Value = 11973.4284
$ws.Range("M122").Value = -5025.571599999999
$ws.Range("N122").Value = -16873.4284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 993.3
$ws.Range("I2").Value = 1643
$ws.Range("J2").Value = 18.75
$ws.Range("K2").Value = 9858
$ws.Range("L2").Value = 112.5
$ws.Range("M2").Value = -9745
$ws.Range("N2").Value = -338.5
$ws.Range("H68").Value = 2448.4565
$ws.Range("J68").Value = 2538.372
$ws.Range("L68").Value = 7615.116
$ws.Range("N68").Value = -9237.116
$ws.Range("H71").Value = 2448.4565
$ws.Range("J71").Value = 2538.372
$ws.Range("L71").Value = 22845.348
$ws.Range("N71").Value = -30957.348
$ws.Range("H132").Value = 2000
$ws.Range("J132").Value = 2000
$ws.Range("L132").Value = 18000
$ws.Range("N132").Value = -23060

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 24482.758
$ws.Range("J47").Value = 10000
$ws.Range("L47").Value = 10000
$ws.Range("N47").Value = -11136
$ws.Range("H97").Value = 683.36365
$ws.Range("I97").Value = 704
$ws.Range("J97").Value = 653.55554
$ws.Range("K97").Value = 704
$ws.Range("L97").Value = 653.55554
$ws.Range("M97").Value = -208
$ws.Range("N97").Value = -1645.55554
$ws.Range("H102").Value = 5658.3335
$ws.Range("I102").Value = 5987.8335
$ws.Range("K102").Value = 5987.8335
$ws.Range("M102").Value = -4365.8335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2224.524
$ws.Range("I46").Value = 1194
$ws.Range("J46").Value = 2636.7334
$ws.Range("K46").Value = 1194
$ws.Range("L46").Value = 2636.7334
$ws.Range("M46").Value = -1006
$ws.Range("N46").Value = -3012.7334
$ws.Range("H61").Value = 9310.166999999999
$ws.Range("I61").Value = 9215.5
$ws.Range("J61").Value = 9499.5
$ws.Range("K61").Value = 9215.5
$ws.Range("L61").Value = 9499.5
$ws.Range("M61").Value = -9013.5
$ws.Range("N61").Value = -9903.5
$ws.Range("H113").Value = 9310.166999999999
$ws.Range("I113").Value = 9215.5
$ws.Range("J113").Value = 9499.5
$ws.Range("K113").Value = 9215.5
$ws.Range("L113").Value = 9499.5
$ws.Range("M113").Value = -7045.5
$ws.Range("N113").Value = -13839.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 40025
$ws.Range("I40").Value = 40025
$ws.Range("K40").Value = 40025
$ws.Range("M40").Value = -39876
$ws.Range("H62").Value = 9500
$ws.Range("I62").Value = 9000
$ws.Range("K62").Value = 9000
$ws.Range("M62").Value = -8376
$ws.Range("H65").Value = 9500
$ws.Range("I65").Value = 9000
$ws.Range("K65").Value = 45000
$ws.Range("M65").Value = -41880
$ws.Range("H132").Value = 6974.857
$ws.Range("I132").Value = 5186.4443
$ws.Range("K132").Value = 15559.3329
$ws.Range("M132").Value = -13029.3329
$ws.Range("H136").Value = 3685.1316
$ws.Range("I136").Value = 3766.7932
$ws.Range("J136").Value = 3422
$ws.Range("K136").Value = 11300.3796
$ws.Range("L136").Value = 10266
$ws.Range("M136").Value = -8750.3796
$ws.Range("N136").Value = -15366

